# The two records currently stored in rows 4 and 5 of the sheet got
# swapped (row 4 now holds what used to be row 5's observation, and vice
# versa). Only the columns whose values actually differ between the two
# records need to be touched; columns that already hold identical values
# in both rows (dates, municipality names, taxon info, etc.) are left
# untouched so we don't risk Excel re-interpreting them (e.g. turning a
# literal date-like text value into a real date serial number).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$swapCols = @("A", "P", "Q", "R", "S", "AW", "AX")

foreach ($col in $swapCols) {
    $cell4 = $ws.Range("$col" + "4")
    $cell5 = $ws.Range("$col" + "5")

    $val4 = $cell4.Value()
    $val5 = $cell5.Value()

    $cell4.Value = $val5
    $cell5.Value = $val4
}
